$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 49

# Leading apostrophe forces the date-looking string to be stored as text
# (matching the original data which stores dates as plain text/inlineStr),
# then reset the style so no extra "quote prefix" formatting sticks to the cell.
$ws.Cells.Item($row, 1).Value = "'2025-10-03"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = 54.45999908447266
$ws.Cells.Item($row, 3).Value = 716.0999755859375
$ws.Cells.Item($row, 4).Value = 328.4500122070312
